$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 117
$ws.Range("A117").Value = 32
$ws.Range("B117").Value = "GT"
$ws.Range("C117").Value = "DC"
$ws.Range("D117").Value = 1
$ws.Range("E117").Value = "GT"
$ws.Range("F117").Value = "DC"
$ws.Range("G117").Value = 5
$ws.Range("H117").Value = "DC"
$ws.Range("I117").Value = "VK Sharma"
$ws.Range("J117").Value = "VKS"
$ws.Range("K117").Value = "Wicket"
$ws.Range("L117").Value = "Not Out"
$ws.Range("M117").Value = "Out"
$ws.Range("N117").Value = "DA Miller"
$ws.Range("O117").Value = "I Sharma"
$ws.Range("P117").Value = "Successful"
$ws.Range("Q117").Value = "No"

# Row 118
$ws.Range("A118").Value = 32
$ws.Range("B118").Value = "GT"
$ws.Range("C118").Value = "DC"
$ws.Range("D118").Value = 1
$ws.Range("E118").Value = "GT"
$ws.Range("F118").Value = "DC"
$ws.Range("G118").Value = 7
$ws.Range("H118").Value = "GT"
$ws.Range("I118").Value = "VK Sharma"
$ws.Range("J118").Value = "VKS"
$ws.Range("K118").Value = "Wide"
$ws.Range("L118").Value = "Not Called"
$ws.Range("M118").Value = "Called"
$ws.Range("N118").Value = "R Tewatia"
$ws.Range("O118").Value = "Mukesh Kumar"
$ws.Range("P118").Value = "Successful"
$ws.Range("Q118").Value = "No"

# Row 119
$ws.Range("A119").Value = 32
$ws.Range("B119").Value = "GT"
$ws.Range("C119").Value = "DC"
$ws.Range("D119").Value = 1
$ws.Range("E119").Value = "GT"
$ws.Range("F119").Value = "DC"
$ws.Range("G119").Value = 12
$ws.Range("H119").Value = "GT"
$ws.Range("I119").Value = "NA Patwardhan"
$ws.Range("J119").Value = "NAP"
$ws.Range("K119").Value = "Wicket"
$ws.Range("L119").Value = "Out"
$ws.Range("M119").Value = "Out"
$ws.Range("N119").Value = "R Tewatia"
$ws.Range("O119").Value = "AR Patel"
$ws.Range("P119").Value = "Unsuccessful"
$ws.Range("Q119").Value = "Yes"

# Row 120
$ws.Range("A120").Value = 32
$ws.Range("B120").Value = "GT"
$ws.Range("C120").Value = "DC"
$ws.Range("D120").Value = 2
$ws.Range("E120").Value = "DC"
$ws.Range("F120").Value = "GT"
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = "GT"
$ws.Range("I120").Value = "NA Patwardhan"
$ws.Range("J120").Value = "NAP"
$ws.Range("K120").Value = "Wicket"
$ws.Range("L120").Value = "Not Out"
$ws.Range("M120").Value = "Not Out"
$ws.Range("N120").Value = "PP Shaw"
$ws.Range("O120").Value = "S Sandeep Warrier"
$ws.Range("P120").Value = "Unsuccessful"
$ws.Range("Q120").Value = "Yes"

# Row 121
$ws.Range("A121").Value = 32
$ws.Range("B121").Value = "GT"
$ws.Range("C121").Value = "DC"
$ws.Range("D121").Value = 2
$ws.Range("E121").Value = "DC"
$ws.Range("F121").Value = "GT"
$ws.Range("G121").Value = 3
$ws.Range("H121").Value = "GT"
$ws.Range("I121").Value = "NA Patwardhan"
$ws.Range("J121").Value = "NAP"
$ws.Range("K121").Value = "Wicket"
$ws.Range("L121").Value = "Not Out"
$ws.Range("M121").Value = "Not Out"
$ws.Range("N121").Value = "PP Shaw"
$ws.Range("O121").Value = "S Sandeep Warrier"
$ws.Range("P121").Value = "Unsuccessful"
$ws.Range("Q121").Value = "No"

# Row 122
$ws.Range("A122").Value = 33
$ws.Range("B122").Value = "PBKS"
$ws.Range("C122").Value = "MI"
$ws.Range("D122").Value = 1
$ws.Range("E122").Value = "MI"
$ws.Range("F122").Value = "PBKS"
$ws.Range("G122").Value = 5
$ws.Range("H122").Value = "MI"
$ws.Range("I122").Value = "VA Kulkarni"
$ws.Range("J122").Value = "VAK"
$ws.Range("K122").Value = "Wicket"
$ws.Range("L122").Value = "Out"
$ws.Range("M122").Value = "Not Out"
$ws.Range("N122").Value = "RG Sharma"
$ws.Range("O122").Value = "HV Patel"
$ws.Range("P122").Value = "Successful"
$ws.Range("Q122").Value = "No"

# Row 123
$ws.Range("A123").Value = 33
$ws.Range("B123").Value = "PBKS"
$ws.Range("C123").Value = "MI"
$ws.Range("D123").Value = 1
$ws.Range("E123").Value = "MI"
$ws.Range("F123").Value = "PBKS"
$ws.Range("G123").Value = 15
$ws.Range("H123").Value = "MI"
$ws.Range("I123").Value = "VA Kulkarni"
$ws.Range("J123").Value = "VAK"
$ws.Range("K123").Value = "Wide"
$ws.Range("L123").Value = "Not Called"
$ws.Range("M123").Value = "Called"
$ws.Range("N123").Value = "SA Yadav"
$ws.Range("O123").Value = "Arshdeep Singh"
$ws.Range("P123").Value = "Successful"
$ws.Range("Q123").Value = "No"

# Row 124
$ws.Range("A124").Value = 33
$ws.Range("B124").Value = "PBKS"
$ws.Range("C124").Value = "MI"
$ws.Range("D124").Value = 1
$ws.Range("E124").Value = "MI"
$ws.Range("F124").Value = "PBKS"
$ws.Range("G124").Value = 16
$ws.Range("H124").Value = "MI"
$ws.Range("I124").Value = "A Nand Kishore"
$ws.Range("J124").Value = "ANK"
$ws.Range("K124").Value = "Wicket"
$ws.Range("L124").Value = "Out"
$ws.Range("M124").Value = "Not Out"
$ws.Range("N124").Value = "SA Yadav"
$ws.Range("O124").Value = "K Rabada"
$ws.Range("P124").Value = "Successful"
$ws.Range("Q124").Value = "No"

# Row 125
$ws.Range("A125").Value = 33
$ws.Range("B125").Value = "PBKS"
$ws.Range("C125").Value = "MI"
$ws.Range("D125").Value = 1
$ws.Range("E125").Value = "MI"
$ws.Range("F125").Value = "PBKS"
$ws.Range("G125").Value = 19
$ws.Range("H125").Value = "MI"
$ws.Range("I125").Value = "VA Kulkarni"
$ws.Range("J125").Value = "VAK"
$ws.Range("K125").Value = "Wide"
$ws.Range("L125").Value = "Not Called"
$ws.Range("M125").Value = "Called"
$ws.Range("N125").Value = "TH David"
$ws.Range("O125").Value = "SM Curran"
$ws.Range("P125").Value = "Successful"
$ws.Range("Q125").Value = "No"

# Row 126
$ws.Range("A126").Value = 33
$ws.Range("B126").Value = "PBKS"
$ws.Range("C126").Value = "MI"
$ws.Range("D126").Value = 1
$ws.Range("E126").Value = "MI"
$ws.Range("F126").Value = "PBKS"
$ws.Range("G126").Value = 19
$ws.Range("H126").Value = "PBKS"
$ws.Range("I126").Value = "VA Kulkarni"
$ws.Range("J126").Value = "VAK"
$ws.Range("K126").Value = "Wide"
$ws.Range("L126").Value = "Called"
$ws.Range("M126").Value = "Not Called"
$ws.Range("N126").Value = "TH David"
$ws.Range("O126").Value = "SM Curran"
$ws.Range("P126").Value = "Successful"
$ws.Range("Q126").Value = "No"

# Row 127
$ws.Range("A127").Value = 33
$ws.Range("B127").Value = "PBKS"
$ws.Range("C127").Value = "MI"
$ws.Range("D127").Value = 1
$ws.Range("E127").Value = "MI"
$ws.Range("F127").Value = "PBKS"
$ws.Range("G127").Value = 20
$ws.Range("H127").Value = "PBKS"
$ws.Range("I127").Value = "A Nand Kishore"
$ws.Range("J127").Value = "ANK"
$ws.Range("K127").Value = "Wicket"
$ws.Range("L127").Value = "Not Out"
$ws.Range("M127").Value = "Not Out"
$ws.Range("N127").Value = "TH David"
$ws.Range("O127").Value = "HV Patel"
$ws.Range("P127").Value = "Unsuccessful"
$ws.Range("Q127").Value = "No"

# Row 128
$ws.Range("A128").Value = 33
$ws.Range("B128").Value = "PBKS"
$ws.Range("C128").Value = "MI"
$ws.Range("D128").Value = 2
$ws.Range("E128").Value = "PBKS"
$ws.Range("F128").Value = "MI"
$ws.Range("G128").Value = 2
$ws.Range("H128").Value = "MI"
$ws.Range("I128").Value = "A Nand Kishore"
$ws.Range("J128").Value = "ANK"
$ws.Range("K128").Value = "Wicket"
$ws.Range("L128").Value = "Not Out"
$ws.Range("M128").Value = "Out"
$ws.Range("N128").Value = "SM Curran"
$ws.Range("O128").Value = "JJ Bumrah"
$ws.Range("P128").Value = "Successful"
$ws.Range("Q128").Value = "No"

# Row 129
$ws.Range("A129").Value = 33
$ws.Range("B129").Value = "PBKS"
$ws.Range("C129").Value = "MI"
$ws.Range("D129").Value = 2
$ws.Range("E129").Value = "PBKS"
$ws.Range("F129").Value = "MI"
$ws.Range("G129").Value = 10
$ws.Range("H129").Value = "PBKS"
$ws.Range("I129").Value = "A Nand Kishore"
$ws.Range("J129").Value = "ANK"
$ws.Range("K129").Value = "Wicket"
$ws.Range("L129").Value = "Out"
$ws.Range("M129").Value = "Out"
$ws.Range("N129").Value = "JM Sharma"
$ws.Range("O129").Value = "A Madhwal"
$ws.Range("P129").Value = "Unsuccessful"
$ws.Range("Q129").Value = "Yes"

# Row 130
$ws.Range("A130").Value = 33
$ws.Range("B130").Value = "PBKS"
$ws.Range("C130").Value = "MI"
$ws.Range("D130").Value = 2
$ws.Range("E130").Value = "PBKS"
$ws.Range("F130").Value = "MI"
$ws.Range("G130").Value = 19
$ws.Range("H130").Value = "PBKS"
$ws.Range("I130").Value = "VA Kulkarni"
$ws.Range("J130").Value = "VAK"
$ws.Range("K130").Value = "Wide"
$ws.Range("L130").Value = "Not Called"
$ws.Range("M130").Value = "Not Called"
$ws.Range("N130").Value = "Harpreet Brar"
$ws.Range("O130").Value = "HH Pandya"
$ws.Range("P130").Value = "Unsuccessful"
$ws.Range("Q130").Value = "No"

$ws.Range("I130").Select()
